$d = $word.ActiveDocument

# --- Short (summary) English program paragraph ---------------------------
# "1) Conceptual basis for transport phenomena study2) General properties of
# fluids 3) Kinematics of fluids.4) Conservation Equations in Integral
# form.5) Differential Equations of Fluid Flow. 6) Boundary Layer Theory.
# 7) Flow in ducts:"
# becomes 7 runs of text separated by manual line breaks.

$d.Content.Find.Execute("phenomena study2) General", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "phenomena study^l2) General", 2) | Out-Null

$d.Content.Find.Execute("of fluids 3) Kinematics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "of fluids ^l3) Kinematics", 2) | Out-Null

$d.Content.Find.Execute("fluids.4) Conservation", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "fluids.^l4) Conservation", 2) | Out-Null

$d.Content.Find.Execute("Integral form.5) Differential", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Integral form.^l5) Differential", 2) | Out-Null

$d.Content.Find.Execute("Fluid Flow. 6) Boundary", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fluid Flow. ^l6) Boundary", 2) | Out-Null

$d.Content.Find.Execute("Layer Theory.7) Flow", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Layer Theory.^l7) Flow", 2) | Out-Null

# --- Long (detailed) English program paragraph ----------------------------
# "1) Conceptual basis for transport phenomena studyFluids and the continuous
# hypothesis. ... formulation.2) General properties of fluids: ...
# compressibility.3) Fluid Kinematics: ... Reynolds number.4) Conservation
# Equations in Integral form: ... Applications.5) Differential Equations of
# Fluid Flow: ... Applications."
# becomes 6 runs of text separated by manual line breaks.

$d.Content.Find.Execute("phenomena studyFluids and", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "phenomena study^lFluids and", 2) | Out-Null

$d.Content.Find.Execute("formulation.2) General properties", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "formulation.^l2) General properties", 2) | Out-Null

$d.Content.Find.Execute("compressibility.3) Fluid Kinematics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "compressibility.^l3) Fluid Kinematics", 2) | Out-Null

$d.Content.Find.Execute("Reynolds number.4) Conservation Equations in Integral form: Flow", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Reynolds number.^l4) Conservation Equations in Integral form: Flow", 2) | Out-Null

$d.Content.Find.Execute("Applications.5) Differential Equations of Fluid Flow: Mass", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Applications.^l5) Differential Equations of Fluid Flow: Mass", 2) | Out-Null

Write-Output "Done"
